# Updates the cryptos list Price (D) and Volume(1h) (E) columns with freshly
# scraped coinranking.com figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.777.45"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "2.398.78"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.63"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +2.57%  "
$ws.Range("D9").Value = "2.404.06"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.17"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "60.524.38"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "2.411.68"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.13"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +7.08%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.17"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.64"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "573.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.11"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.36%  "
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").Value = "0.0₃0940"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.07"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.15"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.30"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.54"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.53%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("E46").Value = "  -4.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.14"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.587"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.36"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.04%  "
